# Apply updated TPM values to the LR-pairs sheet (Vwf-Gp1ba)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 18.975105
$ws.Range("H2").Value = 56.925315
$ws.Range("I2").Value = 0.9552145540969871
$ws.Range("J2").Value = 0.955214554096987
$ws.Range("M2").Value = 0.8229573333333334
$ws.Range("N2").Value = 2.468872
$ws.Range("O2").Value = 0.2440777672676426
$ws.Range("P2").Value = 0.2440777672676426
$ws.Range("Q2").Value = 15.61570181052
$ws.Range("R2").Value = 140.54131629468
$ws.Range("S2").Value = 0.2331466356255494
$ws.Range("T2").Value = 0.2331466356255494

# Row 3
$ws.Range("G3").Value = 18.975105
$ws.Range("H3").Value = 56.925315
$ws.Range("I3").Value = 0.9552145540969871
$ws.Range("J3").Value = 0.955214554096987
$ws.Range("O3").Value = 0.4345811965947162
$ws.Range("P3").Value = 0.4345811965947162
$ws.Range("Q3").Value = 27.80380390419
$ws.Range("R3").Value = 250.23423513771
$ws.Range("S3").Value = 0.4151182839241569
$ws.Range("T3").Value = 0.4151182839241569

# Row 4
$ws.Range("G4").Value = 18.975105
$ws.Range("H4").Value = 56.925315
$ws.Range("I4").Value = 0.9552145540969871
$ws.Range("J4").Value = 0.955214554096987
$ws.Range("M4").Value = 1.083466
$ws.Range("N4").Value = 3.250398
$ws.Range("O4").Value = 0.3213410361376413
$ws.Range("P4").Value = 0.3213410361376413
$ws.Range("Q4").Value = 20.55888111393
$ws.Range("R4").Value = 185.02993002537
$ws.Range("S4").Value = 0.3069496345472808
$ws.Range("T4").Value = 0.3069496345472808

# Row 5
$ws.Range("I5").Value = 0.01570916103663723
$ws.Range("J5").Value = 0.01570916103663723
$ws.Range("M5").Value = 0.8229573333333334
$ws.Range("N5").Value = 2.468872
$ws.Range("O5").Value = 0.2440777672676426
$ws.Range("P5").Value = 0.2440777672676426
$ws.Range("Q5").Value = 0.2568109681635556
$ws.Range("R5").Value = 2.311298713472
$ws.Range("S5").Value = 0.003834256951470261
$ws.Range("T5").Value = 0.003834256951470261

# Row 6
$ws.Range("I6").Value = 0.01570916103663723
$ws.Range("J6").Value = 0.01570916103663723
$ws.Range("O6").Value = 0.4345811965947162
$ws.Range("P6").Value = 0.4345811965947162
$ws.Range("S6").Value = 0.0068269060008009
$ws.Range("T6").Value = 0.006826906000800901

# Row 7
$ws.Range("I7").Value = 0.01570916103663723
$ws.Range("J7").Value = 0.01570916103663723
$ws.Range("M7").Value = 1.083466
$ws.Range("N7").Value = 3.250398
$ws.Range("O7").Value = 0.3213410361376413
$ws.Range("P7").Value = 0.3213410361376413
$ws.Range("Q7").Value = 0.3381049553386666
$ws.Range("R7").Value = 3.042944598048
$ws.Range("S7").Value = 0.00504799808436607
$ws.Range("T7").Value = 0.005047998084366071

# Row 8
$ws.Range("G8").Value = 0.5775933333333333
$ws.Range("H8").Value = 1.73278
$ws.Range("I8").Value = 0.02907628486637583
$ws.Range("J8").Value = 0.02907628486637583
$ws.Range("M8").Value = 0.8229573333333334
$ws.Range("N8").Value = 2.468872
$ws.Range("O8").Value = 0.2440777672676426
$ws.Range("P8").Value = 0.2440777672676426
$ws.Range("Q8").Value = 0.4753346693511111
$ws.Range("R8").Value = 4.278012024160001
$ws.Range("S8").Value = 0.007096874690622959
$ws.Range("T8").Value = 0.007096874690622958

# Row 9
$ws.Range("G9").Value = 0.5775933333333333
$ws.Range("H9").Value = 1.73278
$ws.Range("I9").Value = 0.02907628486637583
$ws.Range("J9").Value = 0.02907628486637583
$ws.Range("O9").Value = 0.4345811965947162
$ws.Range("P9").Value = 0.4345811965947162
$ws.Range("Q9").Value = 0.8463348042799999
$ws.Range("R9").Value = 7.617013238519999
$ws.Range("S9").Value = 0.01263600666975845
$ws.Range("T9").Value = 0.01263600666975845

# Row 10
$ws.Range("G10").Value = 0.5775933333333333
$ws.Range("H10").Value = 1.73278
$ws.Range("I10").Value = 0.02907628486637583
$ws.Range("J10").Value = 0.02907628486637583
$ws.Range("M10").Value = 1.083466
$ws.Range("N10").Value = 3.250398
$ws.Range("O10").Value = 0.3213410361376413
$ws.Range("P10").Value = 0.3213410361376413
$ws.Range("Q10").Value = 0.6258027384933332
$ws.Range("R10").Value = 5.632224646439999
$ws.Range("S10").Value = 0.009343403505994429
$ws.Range("T10").Value = 0.009343403505994429

